{"js": "const replacements = [\n  { row: 0, col: 0, oldText: \"61+14=\", newText: \"15+42=\" },\n  { row: 0, col: 1, oldText: \"79-6=\", newText: \"20+56=\" },\n  { row: 0, col: 2, oldText: \"7+22=\", newText: \"70-47=\" },\n  { row: 0, col: 3, oldText: \"33+26=\", newText: \"32+47=\" },\n  { row: 0, col: 4, oldText: \"0+62=\", newText: \"0+97=\" },\n  { row: 1, col: 0, oldText: \"25+71=\", newText: \"24-16=\" },\n  { row: 1, col: 1, oldText: \"48+15=\", newText: \"15+64=\" },\n  { row: 1, col: 2, oldText: \"36-1=\", newText: \"25+50=\" },\n  { row: 1, col: 3, oldText: \"88-69=\", newText: \"43+31=\" },\n  { row: 1, col: 4, oldText: \"48+7=\", newText: \"20+9=\" },\n  { row: 2, col: 0, oldText: \"23+44=\", newText: \"39+25=\" },\n  { row: 2, col: 1, oldText: \"73-47=\", newText: \"55-29=\" },\n  { row: 2, col: 2, oldText: \"78-9=\", newText: \"66+31=\" },\n  { row: 2, col: 3, oldText: \"75+8=\", newText: \"16+16=\" },\n  { row: 2, col: 4, oldText: \"99-52=\", newText: \"81-72=\" },\n  { row: 3, col: 0, oldText: \"34-25=\", newText: \"18+78=\" },\n  { row: 3, col: 1, oldText: \"0+6=\", newText: \"71-70=\" },\n  { row: 3, col: 2, oldText: \"46+18=\", newText: \"45+38=\" },\n  { row: 3, col: 3, oldText: \"9+78=\", newText: \"54-52=\" },\n  { row: 3, col: 4, oldText: \"46-14=\", newText: \"33+15=\" },\n  { row: 4, col: 0, oldText: \"74+17=\", newText: \"25+9=\" },\n  { row: 4, col: 1, oldText: \"13+23=\", newText: \"9-8=\" },\n  { row: 4, col: 2, oldText: \"13+7=\", newText: \"67-37=\" },\n  { row: 4, col: 3, oldText: \"51-19=\", newText: \"33+33=\" },\n  { row: 4, col: 4, oldText: \"53-4=\", newText: \"73-71=\" },\n  { row: 5, col: 0, oldText: \"75+6=\", newText: \"37+44=\" },\n  { row: 5, col: 1, oldText: \"26-24=\", newText: \"80-10=\" },\n  { row: 5, col: 2, oldText: \"9+51=\", newText: \"0+37=\" },\n  { row: 5, col: 3, oldText: \"92+2=\", newText: \"69-40=\" },\n  { row: 5, col: 4, oldText: \"54+3=\", newText: \"18+57=\" },\n  { row: 6, col: 0, oldText: \"37-13=\", newText: \"12+49=\" },\n  { row: 6, col: 1, oldText: \"6+19=\", newText: \"78-77=\" },\n  { row: 6, col: 2, oldText: \"51+1=\", newText: \"89-46=\" },\n  { row: 6, col: 3, oldText: \"95-90=\", newText: \"56+38=\" },\n  { row: 6, col: 4, oldText: \"38+59=\", newText: \"67+20=\" },\n  { row: 7, col: 0, oldText: \"5+75=\", newText: \"47+51=\" },\n  { row: 7, col: 1, oldText: \"45+16=\", newText: \"21+11=\" },\n  { row: 7, col: 2, oldText: \"48+49=\", newText: \"86-47=\" },\n  { row: 7, col: 3, oldText: \"71-7=\", newText: \"52-25=\" },\n  { row: 7, col: 4, oldText: \"45-15=\", newText: \"77+11=\" },\n  { row: 8, col: 0, oldText: \"33-27=\", newText: \"73-56=\" },\n  { row: 8, col: 1, oldText: \"43-20=\", newText: \"88-67=\" },\n  { row: 8, col: 2, oldText: \"80+2=\", newText: \"84-11=\" },\n  { row: 8, col: 3, oldText: \"67+2=\", newText: \"27-4=\" },\n  { row: 8, col: 4, oldText: \"57+3=\", newText: \"64+0=\" },\n  { row: 9, col: 0, oldText: \"62+30=\", newText: \"70-36=\" },\n  { row: 9, col: 1, oldText: \"69-55=\", newText: \"75-44=\" },\n  { row: 9, col: 2, oldText: \"79-52=\", newText: \"39+55=\" },\n  { row: 9, col: 3, oldText: \"12-9=\", newText: \"91-66=\" },\n  { row: 9, col: 4, oldText: \"12-1=\", newText: \"19+74=\" },\n  { row: 10, col: 0, oldText: \"98-80=\", newText: \"0+51=\" },\n  { row: 10, col: 1, oldText: \"36-5=\", newText: \"13+38=\" },\n  { row: 10, col: 2, oldText: \"75-36=\", newText: \"75-70=\" },\n  { row: 10, col: 3, oldText: \"24-23=\", newText: \"34+3=\" },\n  { row: 10, col: 4, oldText: \"25+17=\", newText: \"7+1=\" },\n  { row: 11, col: 0, oldText: \"67-17=\", newText: \"59-21=\" },\n  { row: 11, col: 1, oldText: \"40+38=\", newText: \"69-7=\" },\n  { row: 11, col: 2, oldText: \"70-1=\", newText: \"0+11=\" },\n  { row: 11, col: 3, oldText: \"1+73=\", newText: \"93-14=\" },\n  { row: 11, col: 4, oldText: \"91-35=\", newText: \"56+37=\" },\n  { row: 12, col: 0, oldText: \"69-4=\", newText: \"38-13=\" },\n  { row: 12, col: 1, oldText: \"9+40=\", newText: \"46+41=\" },\n  { row: 12, col: 2, oldText: \"31+54=\", newText: \"8+49=\" },\n  { row: 12, col: 3, oldText: \"28+1=\", newText: \"9+71=\" },\n  { row: 12, col: 4, oldText: \"23+5=\", newText: \"27-20=\" },\n  { row: 13, col: 0, oldText: \"8+9=\", newText: \"67-56=\" },\n  { row: 13, col: 1, oldText: \"40+3=\", newText: \"47-4=\" },\n  { row: 13, col: 2, oldText: \"36+16=\", newText: \"34+21=\" },\n  { row: 13, col: 3, oldText: \"2+79=\", newText: \"51+24=\" },\n  { row: 13, col: 4, oldText: \"8+80=\", newText: \"90-50=\" },\n  { row: 14, col: 0, oldText: \"35+4=\", newText: \"77-13=\" },\n  { row: 14, col: 1, oldText: \"93-33=\", newText: \"72-12=\" },\n  { row: 14, col: 2, oldText: \"82-80=\", newText: \"53+0=\" },\n  { row: 14, col: 3, oldText: \"59-48=\", newText: \"70-46=\" },\n  { row: 14, col: 4, oldText: \"69-25=\", newText: \"2+45=\" },\n  { row: 15, col: 0, oldText: \"66-56=\", newText: \"24+28=\" },\n  { row: 15, col: 1, oldText: \"74-35=\", newText: \"85-20=\" },\n  { row: 15, col: 2, oldText: \"98-53=\", newText: \"61-59=\" },\n  { row: 15, col: 3, oldText: \"15-13=\", newText: \"89-82=\" },\n  { row: 15, col: 4, oldText: \"43-16=\", newText: \"75+11=\" },\n  { row: 16, col: 0, oldText: \"56+37=\", newText: \"2+75=\" },\n  { row: 16, col: 1, oldText: \"29+43=\", newText: \"91-78=\" },\n  { row: 16, col: 2, oldText: \"63+6=\", newText: \"31+65=\" },\n  { row: 16, col: 3, oldText: \"88-0=\", newText: \"85-32=\" },\n  { row: 16, col: 4, oldText: \"71-20=\", newText: \"55+41=\" },\n  { row: 17, col: 0, oldText: \"42+4=\", newText: \"79+12=\" },\n  { row: 17, col: 1, oldText: \"76-65=\", newText: \"94-68=\" },\n  { row: 17, col: 2, oldText: \"71-28=\", newText: \"83-29=\" },\n  { row: 17, col: 3, oldText: \"21+75=\", newText: \"16+28=\" },\n  { row: 17, col: 4, oldText: \"4-2=\", newText: \"91+3=\" },\n  { row: 18, col: 0, oldText: \"16+60=\", newText: \"44-4=\" },\n  { row: 18, col: 1, oldText: \"69-62=\", newText: \"38-35=\" },\n  { row: 18, col: 2, oldText: \"45-18=\", newText: \"50-11=\" },\n  { row: 18, col: 3, oldText: \"84-0=\", newText: \"32-16=\" },\n  { row: 18, col: 4, oldText: \"28+33=\", newText: \"71+12=\" },\n  { row: 19, col: 0, oldText: \"54-28=\", newText: \"82-75=\" },\n  { row: 19, col: 1, oldText: \"72-25=\", newText: \"26+17=\" },\n  { row: 19, col: 2, oldText: \"55-25=\", newText: \"46-41=\" },\n  { row: 19, col: 3, oldText: \"10+71=\", newText: \"66-57=\" },\n  { row: 19, col: 4, oldText: \"86-11=\", newText: \"74-17=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nlet replacedCount = 0;\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  // Scope the search to this specific cell so identical formulas elsewhere\n  // in the table (or text produced by earlier replacements) are never\n  // mistakenly matched.\n  const found = cell.body.search(r.oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items,text\");\n  await context.sync();\n\n  if (found.items.length > 0) {\n    found.items[0].insertText(r.newText, Word.InsertLocation.replace);\n    replacedCount++;\n  } else {\n    // Fallback: the cell text didn't match what we expected (e.g. it was\n    // already updated) -- just overwrite the whole cell body directly.\n    cell.body.insertText(r.newText, Word.InsertLocation.replace);\n    replacedCount++;\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"61+14=\"; NewText = \"15+42=\" }\n    @{ Row = 1; Col = 2; OldText = \"79-6=\"; NewText = \"20+56=\" }\n    @{ Row = 1; Col = 3; OldText = \"7+22=\"; NewText = \"70-47=\" }\n    @{ Row = 1; Col = 4; OldText = \"33+26=\"; NewText = \"32+47=\" }\n    @{ Row = 1; Col = 5; OldText = \"0+62=\"; NewText = \"0+97=\" }\n    @{ Row = 2; Col = 1; OldText = \"25+71=\"; NewText = \"24-16=\" }\n    @{ Row = 2; Col = 2; OldText = \"48+15=\"; NewText = \"15+64=\" }\n    @{ Row = 2; Col = 3; OldText = \"36-1=\"; NewText = \"25+50=\" }\n    @{ Row = 2; Col = 4; OldText = \"88-69=\"; NewText = \"43+31=\" }\n    @{ Row = 2; Col = 5; OldText = \"48+7=\"; NewText = \"20+9=\" }\n    @{ Row = 3; Col = 1; OldText = \"23+44=\"; NewText = \"39+25=\" }\n    @{ Row = 3; Col = 2; OldText = \"73-47=\"; NewText = \"55-29=\" }\n    @{ Row = 3; Col = 3; OldText = \"78-9=\"; NewText = \"66+31=\" }\n    @{ Row = 3; Col = 4; OldText = \"75+8=\"; NewText = \"16+16=\" }\n    @{ Row = 3; Col = 5; OldText = \"99-52=\"; NewText = \"81-72=\" }\n    @{ Row = 4; Col = 1; OldText = \"34-25=\"; NewText = \"18+78=\" }\n    @{ Row = 4; Col = 2; OldText = \"0+6=\"; NewText = \"71-70=\" }\n    @{ Row = 4; Col = 3; OldText = \"46+18=\"; NewText = \"45+38=\" }\n    @{ Row = 4; Col = 4; OldText = \"9+78=\"; NewText = \"54-52=\" }\n    @{ Row = 4; Col = 5; OldText = \"46-14=\"; NewText = \"33+15=\" }\n    @{ Row = 5; Col = 1; OldText = \"74+17=\"; NewText = \"25+9=\" }\n    @{ Row = 5; Col = 2; OldText = \"13+23=\"; NewText = \"9-8=\" }\n    @{ Row = 5; Col = 3; OldText = \"13+7=\"; NewText = \"67-37=\" }\n    @{ Row = 5; Col = 4; OldText = \"51-19=\"; NewText = \"33+33=\" }\n    @{ Row = 5; Col = 5; OldText = \"53-4=\"; NewText = \"73-71=\" }\n    @{ Row = 6; Col = 1; OldText = \"75+6=\"; NewText = \"37+44=\" }\n    @{ Row = 6; Col = 2; OldText = \"26-24=\"; NewText = \"80-10=\" }\n    @{ Row = 6; Col = 3; OldText = \"9+51=\"; NewText = \"0+37=\" }\n    @{ Row = 6; Col = 4; OldText = \"92+2=\"; NewText = \"69-40=\" }\n    @{ Row = 6; Col = 5; OldText = \"54+3=\"; NewText = \"18+57=\" }\n    @{ Row = 7; Col = 1; OldText = \"37-13=\"; NewText = \"12+49=\" }\n    @{ Row = 7; Col = 2; OldText = \"6+19=\"; NewText = \"78-77=\" }\n    @{ Row = 7; Col = 3; OldText = \"51+1=\"; NewText = \"89-46=\" }\n    @{ Row = 7; Col = 4; OldText = \"95-90=\"; NewText = \"56+38=\" }\n    @{ Row = 7; Col = 5; OldText = \"38+59=\"; NewText = \"67+20=\" }\n    @{ Row = 8; Col = 1; OldText = \"5+75=\"; NewText = \"47+51=\" }\n    @{ Row = 8; Col = 2; OldText = \"45+16=\"; NewText = \"21+11=\" }\n    @{ Row = 8; Col = 3; OldText = \"48+49=\"; NewText = \"86-47=\" }\n    @{ Row = 8; Col = 4; OldText = \"71-7=\"; NewText = \"52-25=\" }\n    @{ Row = 8; Col = 5; OldText = \"45-15=\"; NewText = \"77+11=\" }\n    @{ Row = 9; Col = 1; OldText = \"33-27=\"; NewText = \"73-56=\" }\n    @{ Row = 9; Col = 2; OldText = \"43-20=\"; NewText = \"88-67=\" }\n    @{ Row = 9; Col = 3; OldText = \"80+2=\"; NewText = \"84-11=\" }\n    @{ Row = 9; Col = 4; OldText = \"67+2=\"; NewText = \"27-4=\" }\n    @{ Row = 9; Col = 5; OldText = \"57+3=\"; NewText = \"64+0=\" }\n    @{ Row = 10; Col = 1; OldText = \"62+30=\"; NewText = \"70-36=\" }\n    @{ Row = 10; Col = 2; OldText = \"69-55=\"; NewText = \"75-44=\" }\n    @{ Row = 10; Col = 3; OldText = \"79-52=\"; NewText = \"39+55=\" }\n    @{ Row = 10; Col = 4; OldText = \"12-9=\"; NewText = \"91-66=\" }\n    @{ Row = 10; Col = 5; OldText = \"12-1=\"; NewText = \"19+74=\" }\n    @{ Row = 11; Col = 1; OldText = \"98-80=\"; NewText = \"0+51=\" }\n    @{ Row = 11; Col = 2; OldText = \"36-5=\"; NewText = \"13+38=\" }\n    @{ Row = 11; Col = 3; OldText = \"75-36=\"; NewText = \"75-70=\" }\n    @{ Row = 11; Col = 4; OldText = \"24-23=\"; NewText = \"34+3=\" }\n    @{ Row = 11; Col = 5; OldText = \"25+17=\"; NewText = \"7+1=\" }\n    @{ Row = 12; Col = 1; OldText = \"67-17=\"; NewText = \"59-21=\" }\n    @{ Row = 12; Col = 2; OldText = \"40+38=\"; NewText = \"69-7=\" }\n    @{ Row = 12; Col = 3; OldText = \"70-1=\"; NewText = \"0+11=\" }\n    @{ Row = 12; Col = 4; OldText = \"1+73=\"; NewText = \"93-14=\" }\n    @{ Row = 12; Col = 5; OldText = \"91-35=\"; NewText = \"56+37=\" }\n    @{ Row = 13; Col = 1; OldText = \"69-4=\"; NewText = \"38-13=\" }\n    @{ Row = 13; Col = 2; OldText = \"9+40=\"; NewText = \"46+41=\" }\n    @{ Row = 13; Col = 3; OldText = \"31+54=\"; NewText = \"8+49=\" }\n    @{ Row = 13; Col = 4; OldText = \"28+1=\"; NewText = \"9+71=\" }\n    @{ Row = 13; Col = 5; OldText = \"23+5=\"; NewText = \"27-20=\" }\n    @{ Row = 14; Col = 1; OldText = \"8+9=\"; NewText = \"67-56=\" }\n    @{ Row = 14; Col = 2; OldText = \"40+3=\"; NewText = \"47-4=\" }\n    @{ Row = 14; Col = 3; OldText = \"36+16=\"; NewText = \"34+21=\" }\n    @{ Row = 14; Col = 4; OldText = \"2+79=\"; NewText = \"51+24=\" }\n    @{ Row = 14; Col = 5; OldText = \"8+80=\"; NewText = \"90-50=\" }\n    @{ Row = 15; Col = 1; OldText = \"35+4=\"; NewText = \"77-13=\" }\n    @{ Row = 15; Col = 2; OldText = \"93-33=\"; NewText = \"72-12=\" }\n    @{ Row = 15; Col = 3; OldText = \"82-80=\"; NewText = \"53+0=\" }\n    @{ Row = 15; Col = 4; OldText = \"59-48=\"; NewText = \"70-46=\" }\n    @{ Row = 15; Col = 5; OldText = \"69-25=\"; NewText = \"2+45=\" }\n    @{ Row = 16; Col = 1; OldText = \"66-56=\"; NewText = \"24+28=\" }\n    @{ Row = 16; Col = 2; OldText = \"74-35=\"; NewText = \"85-20=\" }\n    @{ Row = 16; Col = 3; OldText = \"98-53=\"; NewText = \"61-59=\" }\n    @{ Row = 16; Col = 4; OldText = \"15-13=\"; NewText = \"89-82=\" }\n    @{ Row = 16; Col = 5; OldText = \"43-16=\"; NewText = \"75+11=\" }\n    @{ Row = 17; Col = 1; OldText = \"56+37=\"; NewText = \"2+75=\" }\n    @{ Row = 17; Col = 2; OldText = \"29+43=\"; NewText = \"91-78=\" }\n    @{ Row = 17; Col = 3; OldText = \"63+6=\"; NewText = \"31+65=\" }\n    @{ Row = 17; Col = 4; OldText = \"88-0=\"; NewText = \"85-32=\" }\n    @{ Row = 17; Col = 5; OldText = \"71-20=\"; NewText = \"55+41=\" }\n    @{ Row = 18; Col = 1; OldText = \"42+4=\"; NewText = \"79+12=\" }\n    @{ Row = 18; Col = 2; OldText = \"76-65=\"; NewText = \"94-68=\" }\n    @{ Row = 18; Col = 3; OldText = \"71-28=\"; NewText = \"83-29=\" }\n    @{ Row = 18; Col = 4; OldText = \"21+75=\"; NewText = \"16+28=\" }\n    @{ Row = 18; Col = 5; OldText = \"4-2=\"; NewText = \"91+3=\" }\n    @{ Row = 19; Col = 1; OldText = \"16+60=\"; NewText = \"44-4=\" }\n    @{ Row = 19; Col = 2; OldText = \"69-62=\"; NewText = \"38-35=\" }\n    @{ Row = 19; Col = 3; OldText = \"45-18=\"; NewText = \"50-11=\" }\n    @{ Row = 19; Col = 4; OldText = \"84-0=\"; NewText = \"32-16=\" }\n    @{ Row = 19; Col = 5; OldText = \"28+33=\"; NewText = \"71+12=\" }\n    @{ Row = 20; Col = 1; OldText = \"54-28=\"; NewText = \"82-75=\" }\n    @{ Row = 20; Col = 2; OldText = \"72-25=\"; NewText = \"26+17=\" }\n    @{ Row = 20; Col = 3; OldText = \"55-25=\"; NewText = \"46-41=\" }\n    @{ Row = 20; Col = 4; OldText = \"10+71=\"; NewText = \"66-57=\" }\n    @{ Row = 20; Col = 5; OldText = \"86-11=\"; NewText = \"74-17=\" }\n)\n\nforeach ($item in $replacements) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    $r = $cell.Range\n    # A table cell's Range includes the trailing end-of-cell marker; trim it\n    # off so we only touch the visible text and keep the run/paragraph\n    # formatting (font, size, alignment) untouched.\n    $r.MoveEnd(1, -1) | Out-Null\n    if ($r.Text -eq $item.OldText) {\n        $r.Text = $item.NewText\n    } else {\n        # Fallback: cell didn't contain exactly what we expected (e.g.\n        # already updated) -- overwrite it anyway so the target value wins.\n        $r.Text = $item.NewText\n    }\n}\n"}
